# "Support for local worksheet names and constants"
#
# 1. Insert a new worksheet "named_ranges_2" right before the existing
#    "named_ranges" sheet and give it three formula rows that read a
#    local name, a global constant, and another global constant.
# 2. Add a bunch of workbook-level (global) and worksheet-level (local)
#    defined names: CONST_DATE, CONST_FLOAT, CONST_INT, CONST_LOCAL_INT,
#    LOCAL_INT, LOCAL_NAME, out_there_in_the_cold.
# 3. Extend the original "named_ranges" sheet with three more rows that
#    exercise those local/global names.
# 4. Add one more row to "general" (sheet1) showing a number formatted in
#    scientific notation.

$wb = $excel.ActiveWorkbook

# --- 1. new worksheet, inserted immediately before "named_ranges" -------
$namedRanges = $wb.Worksheets.Item("named_ranges")
$namedRanges2 = $wb.Worksheets.Add($namedRanges)
$namedRanges2.Name = "named_ranges_2"

# convenience handles (captured AFTER the sheet insert so they point at
# the right worksheet objects)
$namedRanges = $wb.Worksheets.Item("named_ranges")
$headerError = $wb.Worksheets.Item("header_error")
$general = $wb.Worksheets.Item("general")

# --- 2. defined names -----------------------------------------------------
# global constants
$wb.Names.Add("CONST_DATE", "=43383")
$wb.Names.Add("CONST_FLOAT", "=10.2")
$wb.Names.Add("CONST_INT", "=100")

# CONST_LOCAL_INT: local to named_ranges, then global
$namedRanges.Names.Add("CONST_LOCAL_INT", "=100")
$wb.Names.Add("CONST_LOCAL_INT", "=100")

# LOCAL_INT: local to named_ranges, local to header_error, then global
$namedRanges.Names.Add("LOCAL_INT", "=1000")
$headerError.Names.Add("LOCAL_INT", "=2000")
$wb.Names.Add("LOCAL_INT", "=2000")

# LOCAL_NAME: local to named_ranges, local to header_error, then global
$namedRanges.Names.Add("LOCAL_NAME", '="Hey You"')
$headerError.Names.Add("LOCAL_NAME", '="out there in the cold"')
$wb.Names.Add("LOCAL_NAME", '="out there in the cold"')

# points at the new sheet's A2 cell
$wb.Names.Add("out_there_in_the_cold", '=named_ranges_2!$A$2')

# --- 3. populate named_ranges_2 -------------------------------------------
$namedRanges2.Range("A1").Formula = "=LOCAL_NAME"
$namedRanges2.Range("A2").Formula = "=CONST_INT"
$namedRanges2.Range("A3").Formula = "=LOCAL_INT"

# --- 4. extend named_ranges with rows 8-10 --------------------------------
$namedRanges.Range("A8").Value = "local name"
$namedRanges.Range("B8").Formula = "=CONST_INT"

$namedRanges.Range("A9").Formula = "=LOCAL_NAME"
$namedRanges.Range("B9").Formula = "=CONST_DATE"

$namedRanges.Range("A10").Formula = "=LOCAL_INT"
$namedRanges.Range("B10").Formula = "=CONST_FLOAT"
$namedRanges.Range("B10").NumberFormat = "0.00E+00"

$namedRanges.Range("A9").Select() | Out-Null

# --- 5. extend general (sheet1) with row 7 --------------------------------
$general.Range("A7").Value = "float cient"
$general.Range("B7").Value = -220
$general.Range("B7").NumberFormat = "0.00E+00"

$general.Range("B7").Select() | Out-Null
